$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B4").Value = 3612
$ws.Range("C4").Value = 1335967.5
$ws.Range("D4").Value = 1355388.5
$ws.Range("E4").Value = 1372762.5
$ws.Range("F4").Value = 1388104
$ws.Range("G4").Value = 1403975
$ws.Range("H4").Value = 1416371.5
$ws.Range("I4").Value = 7769340253.07633
$ws.Range("J4").Value = 7736283000.486543
$ws.Range("K4").Value = 7702930757.171782
$ws.Range("L4").Value = 7791774566.322871
$ws.Range("M4").Value = 7706588459.991132
$ws.Range("N4").Value = 8203788654.568119
$ws.Range("O4").Value = 0.7068106312292359
$ws.Range("P4").Value = 0.7397563676633444
$ws.Range("Q4").Value = 0.7015503875968992
$ws.Range("R4").Value = 0.7630121816168328
$ws.Range("S4").Value = 0.5586932447397563
$ws.Range("T4").Value = 5815.515911185213
$ws.Range("U4").Value = 5707.797432608099
$ws.Range("V4").Value = 5611.262514216248
$ws.Range("W4").Value = 5613.24984750629
$ws.Range("X4").Value = 5489.120860407865
$ws.Range("Y4").Value = 5792.116443015211
$ws.Range("Z4").Value = -107.7184785771133
$ws.Range("AA4").Value = -204.2533969689648
$ws.Range("AB4").Value = -202.2660636789224
$ws.Range("AC4").Value = -326.3950507773479
$ws.Range("AD4").Value = -23.39946817000146
$ws.Range("AE4").Value = -0.01852260061225763
$ws.Range("AF4").Value = -0.03512214566830019
$ws.Range("AG4").Value = -0.03478041617767669
$ws.Range("AH4").Value = -0.05612486592110932
$ws.Range("AI4").Value = -0.004023627228840754

$ws.Range("B5").Value = 3725
$ws.Range("C5").Value = 1319801
$ws.Range("D5").Value = 1341202
$ws.Range("E5").Value = 1363421.5
$ws.Range("F5").Value = 1381585.5
$ws.Range("G5").Value = 1397192.5
$ws.Range("H5").Value = 1406038.5
$ws.Range("I5").Value = 8055998951.097355
$ws.Range("J5").Value = 7959841788.538153
$ws.Range("K5").Value = 7948783072.3397
$ws.Range("L5").Value = 8059251459.457801
$ws.Range("M5").Value = 7939627805.884805
$ws.Range("N5").Value = 8407253091.374852
$ws.Range("O5").Value = 0.7157046979865772
$ws.Range("P5").Value = 0.7503355704697987
$ws.Range("Q5").Value = 0.7130201342281879
$ws.Range("R5").Value = 0.7798657718120805
$ws.Range("S5").Value = 0.5892617449664429
$ws.Range("T5").Value = 6103.949725070184
$ws.Range("U5").Value = 5934.856784092294
$ws.Range("V5").Value = 5830.026204178018
$ws.Range("W5").Value = 5833.335294455393
$ws.Range("X5").Value = 5682.558277320272
$ws.Range("Y5").Value = 5979.390387514177
$ws.Range("Z5").Value = -169.0929409778901
$ws.Range("AA5").Value = -273.9235208921664
$ws.Range("AB5").Value = -270.6144306147917
$ws.Range("AC5").Value = -421.391447749912
$ws.Range("AD5").Value = -124.5593375560074
$ws.Range("AE5").Value = -0.02770221718625732
$ws.Range("AF5").Value = -0.0448764379180755
$ws.Range("AG5").Value = -0.04433431512440578
$ws.Range("AH5").Value = -0.0690358647646081
$ws.Range("AI5").Value = -0.02040635050521733

$ws.Range("B6").Value = 3387
$ws.Range("C6").Value = 1279438.5
$ws.Range("D6").Value = 1299649
$ws.Range("E6").Value = 1318799.5
$ws.Range("F6").Value = 1335437
$ws.Range("G6").Value = 1346978.5
$ws.Range("H6").Value = 1350445.5
$ws.Range("I6").Value = 8173875948.201067
$ws.Range("J6").Value = 7986558687.475271
$ws.Range("K6").Value = 7970616591.868944
$ws.Range("L6").Value = 8046481424.204811
$ws.Range("M6").Value = 7949054741.263597
$ws.Range("N6").Value = 8337419851.870849
$ws.Range("O6").Value = 0.7440212577502214
$ws.Range("P6").Value = 0.7676409802184825
$ws.Range("Q6").Value = 0.734573368762917
$ws.Range("R6").Value = 0.7894892235016239
$ws.Range("S6").Value = 0.6259226454089164
$ws.Range("T6").Value = 6388.643102580599
$ws.Range("U6").Value = 6145.1658774602
$ws.Range("V6").Value = 6043.842594624083
$ws.Range("W6").Value = 6025.354564988697
$ws.Range("X6").Value = 5901.396897770526
$ws.Range("Y6").Value = 6173.829193307578
$ws.Range("Z6").Value = -243.4772251203985
$ws.Range("AA6").Value = -344.8005079565155
$ws.Range("AB6").Value = -363.2885375919013
$ws.Range("AC6").Value = -487.2462048100724
$ws.Range("AD6").Value = -214.8139092730207
$ws.Range("AE6").Value = -0.03811094487686273
$ws.Range("AF6").Value = -0.05397085146566383
$ws.Range("AG6").Value = -0.05686474134783903
$ws.Range("AH6").Value = -0.07626755744318481
$ws.Range("AI6").Value = -0.03362434022746552

$ws.Range("B7").Value = 3497
$ws.Range("C7").Value = 1337326
$ws.Range("D7").Value = 1359893
$ws.Range("E7").Value = 1374499.5
$ws.Range("F7").Value = 1390836
$ws.Range("G7").Value = 1402092.5
$ws.Range("H7").Value = 1402158.5
$ws.Range("I7").Value = 9285119972.49873
$ws.Range("J7").Value = 9042625665.008123
$ws.Range("K7").Value = 9012994870.710146
$ws.Range("L7").Value = 9065934735.583771
$ws.Range("M7").Value = 8937067093.97122
$ws.Range("N7").Value = 9306560478.133169
$ws.Range("O7").Value = 0.7489276522733772
$ws.Range("P7").Value = 0.7577923934801258
$ws.Range("Q7").Value = 0.7500714898484415
$ws.Range("R7").Value = 0.8049756934515299
$ws.Range("S7").Value = 0.6454103517300543
$ws.Range("T7").Value = 6943.049019086393
$ws.Range("U7").Value = 6649.512619748851
$ws.Range("V7").Value = 6557.292214882687
$ws.Range("W7").Value = 6518.334825661524
$ws.Range("X7").Value = 6374.092361218122
$ws.Range("Y7").Value = 6637.30988909825
$ws.Range("Z7").Value = -293.5363993375422
$ws.Range("AA7").Value = -385.7568042037055
$ws.Range("AB7").Value = -424.7141934248693
$ws.Range("AC7").Value = -568.9566578682707
$ws.Range("AD7").Value = -305.7391299881428
$ws.Range("AE7").Value = -0.04227773684596103
$ws.Range("AF7").Value = -0.05556014413023191
$ws.Range("AG7").Value = -0.06117113565773957
$ws.Range("AH7").Value = -0.08194622511006511
$ws.Range("AI7").Value = -0.04403528322321626

$ws.Range("B8").Value = 3353
$ws.Range("C8").Value = 1227172
$ws.Range("D8").Value = 1244346.5
$ws.Range("E8").Value = 1253977
$ws.Range("F8").Value = 1262701.5
$ws.Range("G8").Value = 1265795
$ws.Range("H8").Value = 1258617
$ws.Range("I8").Value = 9274138473.364799
$ws.Range("J8").Value = 8979521925.002832
$ws.Range("K8").Value = 8921212200.572533
$ws.Range("L8").Value = 8927661349.425362
$ws.Range("M8").Value = 8769093263.025267
$ws.Range("N8").Value = 9063870295.610519
$ws.Range("O8").Value = 0.7482851178049508
$ws.Range("P8").Value = 0.7488815985684462
$ws.Range("Q8").Value = 0.7363555025350432
$ws.Range("R8").Value = 0.7894422904861318
$ws.Range("S8").Value = 0.6394273784670444
$ws.Range("T8").Value = 7557.325683249617
$ws.Range("U8").Value = 7216.255219107245
$ws.Range("V8").Value = 7114.334792880996
$ws.Range("W8").Value = 7070.286484513847
$ws.Range("X8").Value = 6927.735741589488
$ws.Range("Y8").Value = 7201.452304879498
$ws.Range("Z8").Value = -341.0704641423727
$ws.Range("AA8").Value = -442.990890368621
$ws.Range("AB8").Value = -487.0391987357698
$ws.Range("AC8").Value = -629.5899416601287
$ws.Range("AD8").Value = -355.8733783701191
$ws.Range("AE8").Value = -0.04513110568971979
$ws.Range("AF8").Value = -0.05861741427268186
$ws.Range("AG8").Value = -0.06444597191507362
$ws.Range("AH8").Value = -0.0833085628498953
$ws.Range("AI8").Value = -0.04708985602656934
